$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; this shifts the existing data rows
# (old rows 16-107) down to rows 17-108, carrying their formatting along.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly price-report data.
# (Single-quoted strings are used so literal '$' characters are not
# treated as the start of a variable expansion.)
$ws.Range("A16").Value = 6
$ws.Range("B16").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C16").Value = 'Metropolitana'
$ws.Range("D16").Value = 44847
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 100114007
$ws.Range("G16").Value = 'Jengibre'
$ws.Range("H16").Value = 'Sin especificar'
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14425
$ws.Range("N16").Value = '$/caja 13 kilos'
$ws.Range("O16").Value = 'Perú'
$ws.Range("P16").Value = 1110
$ws.Range("Q16").Value = 13
$ws.Range("R16").Value = 'Hortaliza'
